$d = $word.ActiveDocument

# 1. housing_court_division -> housing_court (3 occurrences throughout doc)
$rng = $d.Content
$found = $rng.Find.Execute("housing_court_division", $false, $false, $false, $false, $false, $true, 1, $false, "housing_court", 2)
Write-Host "housing_court_division replace: $found"

# 2. Wrap "(6) Harm to the public interest:" paragraph with an if/endif jinja block,
#    and merge it with the following paragraph using line breaks instead of a paragraph mark.
$rng2 = $d.Content
$oldText = "(6) Harm to the public interest: " + [char]13 + "{{ public_interest_harm }}."
$newText = "{%p if public_interest_harm_yesno %}" + [char]13 + "(6) Harm to the public interest: " + [char]11 + "{{ public_interest_harm }}." + [char]11 + "{%p endif %}"
$found2 = $rng2.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
Write-Host "public_interest_harm wrap: $found2"

# 3. Merge the record appendix if/endif block into a single paragraph.
$rng3 = $d.Content
$oldText3 = "{%p if record_appendix %}" + [char]13 + "{{ record_appendix.url_for() }}" + [char]13 + "{%p endif %}"
$newText3 = "{{ record_appendix.url_for() }}"
$found3 = $rng3.Find.Execute($oldText3, $false, $false, $false, $false, $false, $true, 1, $false, $newText3, 2)
Write-Host "record_appendix merge: $found3"
